$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skeena")

# --- New row 42 data (added in the order that makes new shared strings
#     land at the expected indices: "Sept 13-19" (122), then the Fulton
#     River note (123, shared with I41), then "FN0944" (124)) ---
$ws.Range("C42").Value = "Sept 13-19"
$ws.Range("C42").NumberFormat = "d-mmm"

$ws.Range("I41").Value = "Fulton River ESSR-sockeye target, selective gear only"
$ws.Range("I42").Value = "Fulton River ESSR-sockeye target, selective gear only"

$ws.Range("A42").Value = "FN0944"
$ws.Range("B42").Value = "Aboriginal"
$ws.Range("D42").Value = "Sockeye"
$ws.Range("E42").Value = "Selective Gear"
$ws.Range("F42").Value = "Region 6-Lake Babine Nation"
$ws.Range("G42").Value = 7

# Leave the selection where the author ended up after adding the new row.
[void]$ws.Activate()
[void]$ws.Range("A43").Select()
